$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "details_bank"
$ws.Range("D1").Value = "details_accNo"
$ws.Range("E1").Value = "details_addr"
$ws.Range("F1").Value = "frnds_name"
$ws.Range("G1").Value = "frnds_best"
